$d = $word.ActiveDocument

function Set-BulletParagraphText($oldText, $newText) {
    # The "What we like" / "What we don't like" bullet paragraphs each
    # contain a leading empty run (<w:r/>) followed by a second run
    # carrying the visible (unformatted) text. A plain Find/Replace on
    # such paragraphs collapses the two indistinguishable, un-formatted
    # runs into a single run, silently dropping the leading empty run.
    # Rebuilding the paragraph body via InsertXML lets us keep that empty
    # run intact while still swapping the visible text, matching the
    # target diff exactly.
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text -eq ($oldText + "`r")) {
            $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
            $xml = '<?xml version="1.0" standalone="yes"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
                '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
                '<w:r/><w:r><w:t>' + $escaped + '</w:t></w:r></w:p>' +
                '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $para.Range.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# Title heading (appears twice in the document - once as the Heading1 at
# the top, once as the bold "meta title" run near the end - both replaced
# identically). wdReplaceAll (the final "2" argument) replaces every
# matching occurrence in the search range in one call.
$d.Content.Find.Execute(
    "Play Dark Vortex Free Slot Game | Yggdrasil Gaming", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Dark Vortex for Free - Exciting Gameplay and Massive Wins", 2)

# "What we like" bullet list
Set-BulletParagraphText "Vortex Reels and Vortex Free Spins features " `
    "More ways to win with stacked wild and high pay symbols"

Set-BulletParagraphText "3,125 ways to win with stacked wilds and high pay symbols" `
    "Vortex Reels and Vortex Free Spins special features"

Set-BulletParagraphText "Win up to 7,318 times your bet on every free spin" `
    "Significant winning potential with 3,125 ways to win"

Set-BulletParagraphText "Exciting graphics and sound effects" `
    "Immersive graphics and theme that enhance the gaming experience"

# "What we don't like" bullet list
Set-BulletParagraphText "No progressive jackpot" `
    "Limited number of paylines compared to other games"

Set-BulletParagraphText "No gamble feature" `
    "Not suitable for players who prefer simpler slot games"

# Italic meta description at the end of document
$d.Content.Find.Execute(
    "Read our review of Dark Vortex, a 5-reels and 243-3,125 paylines slot game packed with unique features and scary theme. Play Dark Vortex free today!",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Dark Vortex for free and experience the thrill of stacked wilds and high pay symbols. Win big with 3,125 ways to win!", 2)
